$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8f898ba839b71fcfd21c8a7097d35288b9c33be/e2e/835265cf-fa44-41d4-b9cf-e1f1f5ba33d7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4eb5180f1959eed6548c0d462e02ce7385e87525/e2e/835265cf-fa44-41d4-b9cf-e1f1f5ba33d7.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4eb5180f1959eed6548c0d462e02ce7385e87525/e2e/835265cf-fa44-41d4-b9cf-e1f1f5ba33d7.md"
$mdDisplay = "835265cf-fa44-41d4-b9cf-e1f1f5ba33d7.md"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Columns.Item(16).ColumnWidth = 39.16666666666667

$ws.Range("I7").Value = $mdDisplay
$ws.Hyperlinks.Add($ws.Range("I7"), $latestMdUrl, "", "", $mdDisplay)
$ws.Range("J7").Value = "835265cf-fa44-41d4-b9cf-e1f1f5ba33d7.89831bbc320f56bcf01f8ed87709a28722bf1bf2.zh-cn.xlf"
$ws.Range("K7").Value = "2016-09-07 17:05:48"
$ws.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Columns.Item(16).ColumnWidth = 39.16666666666667

$ws2.Range("I7").Value = $mdDisplay
$ws2.Hyperlinks.Add($ws2.Range("I7"), $latestMdUrl, "", "", $mdDisplay)
$ws2.Range("J7").Value = "835265cf-fa44-41d4-b9cf-e1f1f5ba33d7.89831bbc320f56bcf01f8ed87709a28722bf1bf2.de-de.xlf"
$ws2.Range("K7").Value = "2016-09-07 17:05:57"
$ws2.Range("P7").Value = $errorDetail

Write-Output "done"
